# Record the test result ("Passed") in column D next to the existing
# Test Case Name / UserName / Password data extracted from the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "Passed"
